# Auto-generated edit script applying numeric updates to Brynhildr_Profits workbook
# Updates currentAveragePrice/LevePrice/LeveProfit columns (H:N) across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1110.5294
$ws.Cells.Item(17, 9).Value = 1051.5454
$ws.Cells.Item(17, 10).Value = 1218.6666
$ws.Cells.Item(17, 11).Value = 3154.6362
$ws.Cells.Item(17, 12).Value = 3655.9998
$ws.Cells.Item(17, 13).Value = -2986.6362
$ws.Cells.Item(17, 14).Value = -3991.9998
$ws.Cells.Item(19, 8).Value = 1724
$ws.Cells.Item(19, 9).Value = 1499
$ws.Cells.Item(19, 11).Value = 1499
$ws.Cells.Item(19, 13).Value = -1324
$ws.Cells.Item(40, 8).Value = 1150.0834
$ws.Cells.Item(40, 9).Value = 1033
$ws.Cells.Item(40, 11).Value = 1033
$ws.Cells.Item(40, 13).Value = -858
$ws.Cells.Item(64, 8).Value = 4100.773
$ws.Cells.Item(64, 10).Value = 3802.8333
$ws.Cells.Item(64, 12).Value = 3802.8333
$ws.Cells.Item(64, 14).Value = -4298.8333
$ws.Cells.Item(67, 8).Value = 4100.773
$ws.Cells.Item(67, 10).Value = 3802.8333
$ws.Cells.Item(67, 12).Value = 3802.8333
$ws.Cells.Item(67, 14).Value = -5518.8333
$ws.Cells.Item(70, 8).Value = 3124.5
$ws.Cells.Item(70, 9).Value = 3499
$ws.Cells.Item(70, 10).Value = 2999.6667
$ws.Cells.Item(70, 11).Value = 10497
$ws.Cells.Item(70, 12).Value = 8999.000100000001
$ws.Cells.Item(70, 13).Value = -10227
$ws.Cells.Item(70, 14).Value = -9539.000100000001
$ws.Cells.Item(73, 8).Value = 3124.5
$ws.Cells.Item(73, 9).Value = 3499
$ws.Cells.Item(73, 10).Value = 2999.6667
$ws.Cells.Item(73, 11).Value = 10497
$ws.Cells.Item(73, 12).Value = 8999.000100000001
$ws.Cells.Item(73, 13).Value = -9561
$ws.Cells.Item(73, 14).Value = -10871.0001
$ws.Cells.Item(86, 8).Value = 24000.4
$ws.Cells.Item(86, 9).Value = 20001
$ws.Cells.Item(86, 10).Value = 26666.666
$ws.Cells.Item(86, 11).Value = 20001
$ws.Cells.Item(86, 12).Value = 26666.666
$ws.Cells.Item(86, 13).Value = -18878
$ws.Cells.Item(86, 14).Value = -28912.666
$ws.Cells.Item(89, 8).Value = 24000.4
$ws.Cells.Item(89, 9).Value = 20001
$ws.Cells.Item(89, 10).Value = 26666.666
$ws.Cells.Item(89, 11).Value = 100005
$ws.Cells.Item(89, 12).Value = 133333.33
$ws.Cells.Item(89, 13).Value = -94389
$ws.Cells.Item(89, 14).Value = -144565.33
$ws.Cells.Item(98, 8).Value = 1781.9
$ws.Cells.Item(98, 9).Value = 1637.7142
$ws.Cells.Item(98, 11).Value = 1637.7142
$ws.Cells.Item(98, 13).Value = -139.7141999999999
$ws.Cells.Item(100, 8).Value = 1650.8235
$ws.Cells.Item(100, 9).Value = 1433.0667
$ws.Cells.Item(100, 11).Value = 1433.0667
$ws.Cells.Item(100, 13).Value = -892.0667000000001
$ws.Cells.Item(113, 8).Value = 3155.4443
$ws.Cells.Item(113, 9).Value = 3112.375
$ws.Cells.Item(113, 11).Value = 3112.375
$ws.Cells.Item(113, 13).Value = 141.625
$ws.Cells.Item(122, 8).Value = 1781.9
$ws.Cells.Item(122, 9).Value = 1637.7142
$ws.Cells.Item(122, 11).Value = 4913.142599999999
$ws.Cells.Item(122, 13).Value = -2463.142599999999
$ws.Cells.Item(138, 8).Value = 3748.975
$ws.Cells.Item(138, 9).Value = 4208.8125
$ws.Cells.Item(138, 11).Value = 12626.4375
$ws.Cells.Item(138, 13).Value = -7486.4375

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(6, 8).Value = 2999.1667
$ws.Cells.Item(11, 8).Value = 4515
$ws.Cells.Item(11, 9).Value = 4515
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = 4515
$ws.Cells.Item(11, 12).Value = 0
$ws.Cells.Item(11, 13).Value = -4371
$ws.Cells.Item(11, 14).ClearContents()
$ws.Cells.Item(74, 8).Value = 5656.34
$ws.Cells.Item(74, 9).Value = 2696.6592
$ws.Cells.Item(74, 11).Value = 2696.6592
$ws.Cells.Item(74, 13).Value = -1822.6592
$ws.Cells.Item(77, 8).Value = 5656.34
$ws.Cells.Item(77, 9).Value = 2696.6592
$ws.Cells.Item(77, 11).Value = 13483.296
$ws.Cells.Item(77, 13).Value = -9115.296
$ws.Cells.Item(109, 8).Value = 69500
$ws.Cells.Item(109, 10).Value = 69500
$ws.Cells.Item(109, 12).Value = 69500
$ws.Cells.Item(109, 14).Value = -72274
$ws.Cells.Item(110, 8).Value = 1049.5217
$ws.Cells.Item(110, 9).Value = 1038.7142
$ws.Cells.Item(110, 11).Value = 1038.7142
$ws.Cells.Item(110, 13).Value = 1006.2858
$ws.Cells.Item(132, 8).Value = 5434.891
$ws.Cells.Item(132, 9).Value = 3597.625
$ws.Cells.Item(132, 11).Value = 10792.875
$ws.Cells.Item(132, 13).Value = -8262.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(7, 8).Value = 2000765
$ws.Cells.Item(7, 9).Value = 3000000
$ws.Cells.Item(7, 10).Value = 2295
$ws.Cells.Item(7, 11).Value = 3000000
$ws.Cells.Item(7, 12).Value = 2295
$ws.Cells.Item(7, 13).Value = -2999887
$ws.Cells.Item(7, 14).Value = -2521
$ws.Cells.Item(8, 8).Value = 100
$ws.Cells.Item(8, 9).Value = 100
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = 100
$ws.Cells.Item(8, 12).Value = 0
$ws.Cells.Item(8, 13).Value = 40
$ws.Cells.Item(8, 14).ClearContents()
$ws.Cells.Item(10, 8).Value = 4200
$ws.Cells.Item(10, 9).Value = 1900
$ws.Cells.Item(10, 10).Value = 6500
$ws.Cells.Item(10, 11).Value = 1900
$ws.Cells.Item(10, 12).Value = 6500
$ws.Cells.Item(10, 13).Value = -1760
$ws.Cells.Item(10, 14).Value = -6780
$ws.Cells.Item(12, 8).Value = 275
$ws.Cells.Item(12, 9).Value = 275
$ws.Cells.Item(12, 11).Value = 275
$ws.Cells.Item(12, 13).Value = -107
$ws.Cells.Item(17, 8).Value = 2147.5
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 2147.5
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).Value = 2147.5
$ws.Cells.Item(17, 13).ClearContents()
$ws.Cells.Item(17, 14).Value = -2491.5
$ws.Cells.Item(20, 8).Value = 53611.2
$ws.Cells.Item(20, 9).Value = 95365.09
$ws.Cells.Item(20, 10).Value = 2578.6667
$ws.Cells.Item(20, 11).Value = 95365.09
$ws.Cells.Item(20, 12).Value = 2578.6667
$ws.Cells.Item(20, 13).Value = -95118.09
$ws.Cells.Item(20, 14).Value = -3072.6667
$ws.Cells.Item(128, 8).Value = 7281
$ws.Cells.Item(128, 9).Value = 7281
$ws.Cells.Item(128, 11).Value = 21843
$ws.Cells.Item(128, 13).Value = -19353

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 10200.308
$ws.Cells.Item(16, 9).Value = 10200.308
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 10200.308
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).Value = -9913.308000000001
$ws.Cells.Item(16, 14).ClearContents()
$ws.Cells.Item(18, 8).Value = 51332.668
$ws.Cells.Item(18, 9).Value = 0
$ws.Cells.Item(18, 11).Value = 0
$ws.Cells.Item(18, 13).ClearContents()
$ws.Cells.Item(57, 8).Value = 26000
$ws.Cells.Item(57, 10).Value = 26000
$ws.Cells.Item(57, 12).Value = 26000
$ws.Cells.Item(57, 14).Value = -27120
$ws.Cells.Item(58, 8).Value = 6647.5
$ws.Cells.Item(58, 9).Value = 4615.0835
$ws.Cells.Item(58, 11).Value = 4615.0835
$ws.Cells.Item(58, 13).Value = -4412.0835
$ws.Cells.Item(86, 8).Value = 223846.75
$ws.Cells.Item(86, 10).Value = 3249.5
$ws.Cells.Item(86, 12).Value = 3249.5
$ws.Cells.Item(86, 14).Value = -5495.5
$ws.Cells.Item(89, 8).Value = 223846.75
$ws.Cells.Item(89, 10).Value = 3249.5
$ws.Cells.Item(89, 12).Value = 16247.5
$ws.Cells.Item(89, 14).Value = -27479.5
$ws.Cells.Item(105, 8).Value = 9566.154
$ws.Cells.Item(105, 9).Value = 10196.667
$ws.Cells.Item(105, 10).Value = 2000
$ws.Cells.Item(105, 11).Value = 10196.667
$ws.Cells.Item(105, 12).Value = 2000
$ws.Cells.Item(105, 13).Value = -8449.666999999999
$ws.Cells.Item(105, 14).Value = -5494
$ws.Cells.Item(113, 8).Value = 10200.308
$ws.Cells.Item(113, 9).Value = 10200.308
$ws.Cells.Item(113, 10).Value = 0
$ws.Cells.Item(113, 11).Value = 10200.308
$ws.Cells.Item(113, 12).Value = 0
$ws.Cells.Item(113, 13).Value = -8030.308000000001
$ws.Cells.Item(113, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 9560.182000000001
$ws.Cells.Item(122, 9).Value = 2018.037
$ws.Cells.Item(122, 10).Value = 43499.832
$ws.Cells.Item(122, 11).Value = 6054.111
$ws.Cells.Item(122, 12).Value = 130499.496
$ws.Cells.Item(122, 13).Value = -3604.111
$ws.Cells.Item(122, 14).Value = -135399.496
$ws.Cells.Item(132, 8).Value = 3436
$ws.Cells.Item(132, 9).Value = 3436
$ws.Cells.Item(132, 11).Value = 10308
$ws.Cells.Item(132, 13).Value = -7778
$ws.Cells.Item(136, 8).Value = 6647.5
$ws.Cells.Item(136, 9).Value = 4615.0835
$ws.Cells.Item(136, 11).Value = 13845.2505
$ws.Cells.Item(136, 13).Value = -11295.2505

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1537.0769
$ws.Cells.Item(5, 9).Value = 1204.7142
$ws.Cells.Item(5, 10).Value = 1924.8334
$ws.Cells.Item(5, 11).Value = 3614.1426
$ws.Cells.Item(5, 12).Value = 5774.5002
$ws.Cells.Item(5, 13).Value = -3502.1426
$ws.Cells.Item(5, 14).Value = -5998.5002
$ws.Cells.Item(25, 8).Value = 250125.5
$ws.Cells.Item(25, 9).Value = 500
$ws.Cells.Item(25, 11).Value = 1500
$ws.Cells.Item(25, 13).Value = -1331
$ws.Cells.Item(30, 8).Value = 250125.5
$ws.Cells.Item(30, 9).Value = 500
$ws.Cells.Item(30, 11).Value = 1500
$ws.Cells.Item(30, 13).Value = -1398
$ws.Cells.Item(69, 8).Value = 3450
$ws.Cells.Item(69, 9).Value = 1900
$ws.Cells.Item(69, 10).Value = 5000
$ws.Cells.Item(69, 11).Value = 5700
$ws.Cells.Item(69, 12).Value = 15000
$ws.Cells.Item(69, 13).Value = -4889
$ws.Cells.Item(69, 14).Value = -16622
$ws.Cells.Item(72, 8).Value = 3450
$ws.Cells.Item(72, 9).Value = 1900
$ws.Cells.Item(72, 10).Value = 5000
$ws.Cells.Item(72, 11).Value = 17100
$ws.Cells.Item(72, 12).Value = 45000
$ws.Cells.Item(72, 13).Value = -13044
$ws.Cells.Item(72, 14).Value = -53112
$ws.Cells.Item(74, 8).Value = 14832.167
$ws.Cells.Item(74, 9).Value = 13993
$ws.Cells.Item(74, 11).Value = 41979
$ws.Cells.Item(74, 13).Value = -40918
$ws.Cells.Item(77, 8).Value = 14832.167
$ws.Cells.Item(77, 9).Value = 13993
$ws.Cells.Item(77, 11).Value = 125937
$ws.Cells.Item(77, 13).Value = -120633
$ws.Cells.Item(88, 8).Value = 52666
$ws.Cells.Item(88, 9).Value = 100000
$ws.Cells.Item(88, 10).Value = 28999
$ws.Cells.Item(88, 11).Value = 300000
$ws.Cells.Item(88, 12).Value = 86997
$ws.Cells.Item(88, 13).Value = -299572
$ws.Cells.Item(88, 14).Value = -87853
$ws.Cells.Item(91, 8).Value = 52666
$ws.Cells.Item(91, 9).Value = 100000
$ws.Cells.Item(91, 10).Value = 28999
$ws.Cells.Item(91, 11).Value = 300000
$ws.Cells.Item(91, 12).Value = 86997
$ws.Cells.Item(91, 13).Value = -298518
$ws.Cells.Item(91, 14).Value = -89961
$ws.Cells.Item(113, 8).Value = 790.2174
$ws.Cells.Item(113, 9).Value = 412
$ws.Cells.Item(113, 10).Value = 869.8421
$ws.Cells.Item(113, 11).Value = 1236
$ws.Cells.Item(113, 12).Value = 2609.5263
$ws.Cells.Item(113, 13).Value = 934
$ws.Cells.Item(113, 14).Value = -6949.5263
$ws.Cells.Item(131, 8).Value = 1994.3214
$ws.Cells.Item(131, 9).Value = 759.625
$ws.Cells.Item(131, 10).Value = 2124.2896
$ws.Cells.Item(131, 11).Value = 2278.875
$ws.Cells.Item(131, 12).Value = 6372.8688
$ws.Cells.Item(131, 13).Value = 2761.125
$ws.Cells.Item(131, 14).Value = -16452.8688
$ws.Cells.Item(135, 8).Value = 1537.0769
$ws.Cells.Item(135, 9).Value = 1204.7142
$ws.Cells.Item(135, 10).Value = 1924.8334
$ws.Cells.Item(135, 11).Value = 10842.4278
$ws.Cells.Item(135, 12).Value = 17323.5006
$ws.Cells.Item(135, 13).Value = -8307.427799999999
$ws.Cells.Item(135, 14).Value = -22393.5006
$ws.Cells.Item(139, 8).Value = 10220.5625
$ws.Cells.Item(139, 9).Value = 8615
$ws.Cells.Item(139, 10).Value = 10449.929
$ws.Cells.Item(139, 11).Value = 25845
$ws.Cells.Item(139, 12).Value = 31349.787
$ws.Cells.Item(139, 13).Value = -20705
$ws.Cells.Item(139, 14).Value = -41629.787
$ws.Cells.Item(140, 8).Value = 2299.8
$ws.Cells.Item(140, 9).Value = 1874.75
$ws.Cells.Item(140, 11).Value = 5624.25
$ws.Cells.Item(140, 13).Value = -444.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(36, 8).Value = 2000
$ws.Cells.Item(36, 9).Value = 2000
$ws.Cells.Item(36, 10).Value = 0
$ws.Cells.Item(36, 11).Value = 2000
$ws.Cells.Item(36, 12).Value = 0
$ws.Cells.Item(36, 13).Value = -1515
$ws.Cells.Item(36, 14).ClearContents()
$ws.Cells.Item(80, 8).Value = 3348.6
$ws.Cells.Item(80, 9).Value = 3185.75
$ws.Cells.Item(80, 10).Value = 4000
$ws.Cells.Item(80, 11).Value = 3185.75
$ws.Cells.Item(80, 12).Value = 4000
$ws.Cells.Item(80, 13).Value = -2187.75
$ws.Cells.Item(80, 14).Value = -5996
$ws.Cells.Item(83, 8).Value = 3348.6
$ws.Cells.Item(83, 9).Value = 3185.75
$ws.Cells.Item(83, 10).Value = 4000
$ws.Cells.Item(83, 11).Value = 15928.75
$ws.Cells.Item(83, 12).Value = 20000
$ws.Cells.Item(83, 13).Value = -10936.75
$ws.Cells.Item(83, 14).Value = -29984
$ws.Cells.Item(102, 8).Value = 2397.2415
$ws.Cells.Item(102, 9).Value = 2710.9092
$ws.Cells.Item(102, 11).Value = 2710.9092
$ws.Cells.Item(102, 13).Value = -1088.9092
$ws.Cells.Item(113, 8).Value = 1751.8422
$ws.Cells.Item(113, 9).Value = 1734.4117
$ws.Cells.Item(113, 10).Value = 1900
$ws.Cells.Item(113, 11).Value = 1734.4117
$ws.Cells.Item(113, 12).Value = 1900
$ws.Cells.Item(113, 13).Value = 435.5882999999999
$ws.Cells.Item(113, 14).Value = -6240
$ws.Cells.Item(132, 8).Value = 15262.3
$ws.Cells.Item(132, 9).Value = 22587.334
$ws.Cells.Item(132, 11).Value = 67762.00199999999
$ws.Cells.Item(132, 13).Value = -65232.00199999999
$ws.Cells.Item(134, 8).Value = 44999
$ws.Cells.Item(134, 10).Value = 44999
$ws.Cells.Item(134, 12).Value = 134997
$ws.Cells.Item(134, 14).Value = -140067

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 7001.6
$ws.Cells.Item(7, 9).Value = 6667.6665
$ws.Cells.Item(7, 10).Value = 7502.5
$ws.Cells.Item(7, 11).Value = 6667.6665
$ws.Cells.Item(7, 12).Value = 7502.5
$ws.Cells.Item(7, 13).Value = -6555.6665
$ws.Cells.Item(7, 14).Value = -7726.5
$ws.Cells.Item(22, 8).Value = 3000.8286
$ws.Cells.Item(22, 9).Value = 2357.3125
$ws.Cells.Item(22, 10).Value = 3542.7368
$ws.Cells.Item(22, 11).Value = 2357.3125
$ws.Cells.Item(22, 12).Value = 3542.7368
$ws.Cells.Item(22, 13).Value = -2062.3125
$ws.Cells.Item(22, 14).Value = -4132.736800000001
$ws.Cells.Item(27, 8).Value = 3000.8286
$ws.Cells.Item(27, 9).Value = 2357.3125
$ws.Cells.Item(27, 10).Value = 3542.7368
$ws.Cells.Item(27, 11).Value = 2357.3125
$ws.Cells.Item(27, 12).Value = 3542.7368
$ws.Cells.Item(27, 13).Value = -2250.3125
$ws.Cells.Item(27, 14).Value = -3756.7368
$ws.Cells.Item(46, 8).Value = 4423.3706
$ws.Cells.Item(46, 9).Value = 2734.5
$ws.Cells.Item(46, 10).Value = 4558.48
$ws.Cells.Item(46, 11).Value = 2734.5
$ws.Cells.Item(46, 12).Value = 4558.48
$ws.Cells.Item(46, 13).Value = -2546.5
$ws.Cells.Item(46, 14).Value = -4934.48
$ws.Cells.Item(55, 8).Value = 1064.6809
$ws.Cells.Item(55, 10).Value = 1187.3846
$ws.Cells.Item(55, 12).Value = 1187.3846
$ws.Cells.Item(55, 14).Value = -1533.3846
$ws.Cells.Item(61, 8).Value = 6899.186
$ws.Cells.Item(61, 9).Value = 6833.472
$ws.Cells.Item(61, 10).Value = 7237.143
$ws.Cells.Item(61, 11).Value = 6833.472
$ws.Cells.Item(61, 12).Value = 7237.143
$ws.Cells.Item(61, 13).Value = -6631.472
$ws.Cells.Item(61, 14).Value = -7641.143
$ws.Cells.Item(82, 8).Value = 3684.75
$ws.Cells.Item(82, 10).Value = 0
$ws.Cells.Item(82, 12).Value = 0
$ws.Cells.Item(82, 14).ClearContents()
$ws.Cells.Item(85, 8).Value = 3684.75
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 14).ClearContents()
$ws.Cells.Item(113, 8).Value = 6899.186
$ws.Cells.Item(113, 9).Value = 6833.472
$ws.Cells.Item(113, 10).Value = 7237.143
$ws.Cells.Item(113, 11).Value = 6833.472
$ws.Cells.Item(113, 12).Value = 7237.143
$ws.Cells.Item(113, 13).Value = -4663.472
$ws.Cells.Item(113, 14).Value = -11577.143
$ws.Cells.Item(122, 8).Value = 4797.923
$ws.Cells.Item(122, 9).Value = 3403.7778
$ws.Cells.Item(122, 11).Value = 10211.3334
$ws.Cells.Item(122, 13).Value = -7761.3334
$ws.Cells.Item(126, 8).Value = 7001.6
$ws.Cells.Item(126, 9).Value = 6667.6665
$ws.Cells.Item(126, 10).Value = 7502.5
$ws.Cells.Item(126, 11).Value = 20002.9995
$ws.Cells.Item(126, 12).Value = 22507.5
$ws.Cells.Item(126, 13).Value = -17532.9995
$ws.Cells.Item(126, 14).Value = -27447.5
$ws.Cells.Item(132, 8).Value = 3096.55
$ws.Cells.Item(132, 9).Value = 2415
$ws.Cells.Item(132, 11).Value = 7245
$ws.Cells.Item(132, 13).Value = -4715

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(29, 8).Value = 17495
$ws.Cells.Item(29, 10).Value = 19990
$ws.Cells.Item(29, 12).Value = 19990
$ws.Cells.Item(29, 14).Value = -20570
$ws.Cells.Item(122, 8).Value = 50071.418
$ws.Cells.Item(122, 9).Value = 3284.1052
$ws.Cells.Item(122, 10).Value = 227863.2
$ws.Cells.Item(122, 11).Value = 9852.3156
$ws.Cells.Item(122, 12).Value = 683589.6000000001
$ws.Cells.Item(122, 13).Value = -7402.3156
$ws.Cells.Item(122, 14).Value = -688489.6000000001
